# Appetite - Test Cases - Sprint 7
# Commit: Test case for ID #156674345 - appetite icon and splash screen validation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint-5")

# --- Summary counters (Number of test cases / Test coverage) ---
$ws.Range("B3").Value = 2
$ws.Range("D4").Value = 2

# --- TC1 (row 8): Icon check ---
$ws.Range("B8").Value = "Icon: check if icon is displayed"
$ws.Range("C8").Value = "Check if appetite icon appears in the device's screen"
$ws.Range("D8").Value = "Appetite icon (a fork with a small tomato in a red background) appears in the device's screen."
$ws.Rows.Item(8).RowHeight = 13

# --- TC2 (row 9): Splash screen animation ---
$ws.Range("B9").Value = "Animation: Splash screen"
$ws.Range("C9").Value = "Push the appetite icon"
$ws.Range("D9").Value = "Splash screen shows appetite's logo raising from the bottom."
$ws.Rows.Item(9).RowHeight = 13

# --- Test date text updated globally (shared string) ---
$ws.Range("H8").Value = "04/17/2018"
$ws.Range("H9").Value = "04/17/2018"

# --- TC3 (row 10): cleared out ---
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = ""
$ws.Rows.Item(10).RowHeight = 13

# --- TC4 (row 11): cleared out ---
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""
$ws.Rows.Item(11).RowHeight = 13

# --- TC5 (row 12): cleared out ---
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""
$ws.Rows.Item(12).AutoFit()

# --- View: scroll right, select I10:I12 ---
$ws.Activate()
$ws.Range("I10:I12").Select()

# --- Workbook window position (best effort; host-app window metrics) ---
$excel.ActiveWindow.Left = 0
